$d = $word.ActiveDocument

# Locate the paragraph that ends with "LOB1225: ... (Requisito fraco)"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "LOB1225:.*Requisito fraco") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    $startPara = $d.Paragraphs.Item($targetIndex + 1)
    # The three trailing paragraphs to remove: blank line, "Ver no Jupiter..." line,
    # and the "(c) 2020 ... Creative Commons Attribution" line.
    $endPara = $d.Paragraphs.Item($targetIndex + 3)

    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
